$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.129.14"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.80%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.567.25"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.01%  "
$ws.Range("E4").Value = "  +0.75%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.94%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.490"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.49%  "
$ws.Range("E7").Value = "  +0.87%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "21.98"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.22%  "
$ws.Range("E9").Value = "  +0.25%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0862"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.67%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.789.29"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.91%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.566.76"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.20%  "
$ws.Range("E14").Value = "  +0.41%  "
$ws.Range("E15").Value = "  -0.05%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.130.34"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.84%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.00"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.64%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0703"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "214.61"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.37"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.97%  "
$ws.Range("E21").Value = "  +0.88%  "
$ws.Range("E22").Value = "  +1.19%  "
$ws.Range("E23").Value = "  +0.17%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.32"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.57%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.59"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.04"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.42%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.105"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.02%  "
$ws.Range("E29").Value = "  +0.76%  "
$ws.Range("E30").Value = "  +5.08%  "
$ws.Range("E31").Value = "  +0.52%  "
$ws.Range("E32").Value = "  +0.76%  "
$ws.Range("E33").Value = "  +2.86%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.431.43"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.44%  "
$ws.Range("E35").Value = "  +11.60%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.60"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.54%  "
$ws.Range("E37").Value = "  +2.29%  "
$ws.Range("E38").Value = "  +0.93%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.530"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.42%  "
$ws.Range("E40").Value = "  +3.17%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.807"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.07%  "
$ws.Range("E42").Value = "  +0.81%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.35"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.02%  "
$ws.Range("E44").Value = "  +0.50%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.51"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("E46").Value = "  +0.77%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.706.64"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.13%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "85.92"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.39%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₆0102"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.37%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0520"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.13%  "
$ws.Range("E51").Value = "  -0.14%  "
